# Update "想去人数" (want-to-go count) figures in column F across sheets,
# reflecting refreshed data pulled at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibition) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 5623
$ws1.Range("F6").Value = 936
$ws1.Range("F10").Value = 146
$ws1.Range("F14").Value = 2362
$ws1.Range("F15").Value = 350

# --- Sheet "演出" (Show) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 104

# --- Sheet "全部类型" (All Types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 5623
$ws4.Range("F4").Value = 104
$ws4.Range("F8").Value = 936
$ws4.Range("F12").Value = 146
$ws4.Range("F17").Value = 2362
$ws4.Range("F18").Value = 350
